$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.925.37"
$ws.Range("E2").Value = "  -3.54%  "
$ws.Range("D3").Value = "'2.287.50"
$ws.Range("E3").Value = "  -4.09%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'533.31"
$ws.Range("E5").Value = "  -4.36%  "
$ws.Range("D6").Value = "'130.57"
$ws.Range("E6").Value = "  -2.65%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.580"
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("D9").Value = "'2.287.69"
$ws.Range("E9").Value = "  -3.93%  "
$ws.Range("E10").Value = "  -6.47%  "
$ws.Range("E11").Value = "  -4.51%  "
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("E13").Value = "  -4.26%  "
$ws.Range("D14").Value = "'23.40"
$ws.Range("E14").Value = "  -4.56%  "
$ws.Range("D15").Value = "'2.695.65"
$ws.Range("E15").Value = "  -4.05%  "
$ws.Range("D16").Value = "'57.863.83"
$ws.Range("E16").Value = "  -3.53%  "
$ws.Range("E17").Value = "  -5.17%  "
$ws.Range("D18").Value = "'2.291.39"
$ws.Range("E18").Value = "  -3.96%  "
$ws.Range("E19").Value = "  -5.77%  "
$ws.Range("E20").Value = "  -6.26%  "
$ws.Range("D21").Value = "'312.03"
$ws.Range("E21").Value = "  -3.06%  "
$ws.Range("E22").Value = "  -4.72%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "'62.42"
$ws.Range("E24").Value = "  -2.66%  "
$ws.Range("E25").Value = "  -3.85%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "'7.98"
$ws.Range("E27").Value = "  -5.71%  "
$ws.Range("E28").Value = "  -7.10%  "
$ws.Range("D29").Value = "'170.94"
$ws.Range("E29").Value = "  +0.51%  "
$ws.Range("E30").Value = "  -6.18%  "
$ws.Range("E31").Value = "  -6.15%  "
$ws.Range("E32").Value = "  -5.90%  "
$ws.Range("E33").Value = "  -7.03%  "
$ws.Range("E34").Value = "  -5.49%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  -2.67%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E38").Value = "  -7.88%  "
$ws.Range("E39").Value = "  -6.69%  "
$ws.Range("D40").Value = "'38.13"
$ws.Range("E40").Value = "  -1.36%  "
$ws.Range("E41").Value = "  -7.01%  "
$ws.Range("D42").Value = "'141.33"
$ws.Range("E42").Value = "  -2.75%  "
$ws.Range("D43").Value = "'286.65"
$ws.Range("E43").Value = "  -10.47%  "
$ws.Range("D44").Value = "'3.40"
$ws.Range("E44").Value = "  -4.07%  "
$ws.Range("D45").Value = "'0.0945"
$ws.Range("E45").Value = "  -2.75%  "
$ws.Range("E46").Value = "  -3.40%  "
$ws.Range("E47").Value = "  -2.96%  "
$ws.Range("D48").Value = "'18.05"
$ws.Range("E48").Value = "  -9.08%  "
$ws.Range("D49").Value = "'0.0210"
$ws.Range("E49").Value = "  -4.17%  "
$ws.Range("D50").Value = "'10.94"
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("D51").Value = "'0.0₆0201"
$ws.Range("E51").Value = "  +84.60%  "
